$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$preprocess = 'convert to lower, convert unicode to ascii, remove multiple spaces, trim "space" and ","'
$modelDetails = '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 300'
$features = '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit'
$model = 'Neuron Network'
$templateFilter = '0 filters: '

$rows = @(
    @{ Row = 7;  Time = '20160405_102700'; RunningTime = 253.297; G = 0.979333333333333; H = 0.996699669966997; J = 0.13265306122449 },
    @{ Row = 8;  Time = '20160405_103113'; RunningTime = 249.092; G = 0.974666666666667; H = 0.996699669966997; J = 0.142857142857143 },
    @{ Row = 9;  Time = '20160405_103522'; RunningTime = 273.641; G = 0.971333333333333; H = 0.996699669966997; J = 0.183673469387755 },
    @{ Row = 10; Time = '20160405_103956'; RunningTime = 368.064; G = 0.977333333333333; H = 0.996699669966997; J = 0.122448979591837 },
    @{ Row = 11; Time = '20160405_104604'; RunningTime = 494.456; G = 0.974;              H = 0.996699669966997; J = 0.153061224489796 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Time
    $ws.Range("B$row").Value = $r.RunningTime
    $ws.Range("C$row").Value = $preprocess
    $ws.Range("D$row").Value = $features
    $ws.Range("E$row").Value = $model
    $ws.Range("F$row").Value = $modelDetails
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $templateFilter
    $ws.Range("J$row").Value = $r.J
}
